$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300639748573303
$ws.Range("B1").Value = 1.820849657058716
$ws.Range("C1").Value = 1.722445607185364
$ws.Range("D1").Value = 4.992780208587646
$ws.Range("E1").Value = 1.348918080329895
